$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 110, pushing existing rows 110-188 down to 111-189
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new record
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 44651
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = "Fruta"
$ws.Cells.Item(110, 7).Value = 100102
$ws.Cells.Item(110, 8).Value = "Cítricos"
$ws.Cells.Item(110, 9).Value = 100102004
$ws.Cells.Item(110, 10).Value = "Mandarina"
$ws.Cells.Item(110, 11).Value = "Murcott"
$ws.Cells.Item(110, 12).Value = "Primera"
$ws.Cells.Item(110, 13).Value = 400
$ws.Cells.Item(110, 14).Value = 12500
$ws.Cells.Item(110, 15).Value = 13000
$ws.Cells.Item(110, 16).Value = 12750
$ws.Cells.Item(110, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(110, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(110, 19).Value = 1275
$ws.Cells.Item(110, 20).Value = 10
